$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.3958491905232238
$ws.Range("C2").Value = -0.217077417213067
$ws.Range("D2").Value = -0.1238165243142643
$ws.Range("E2").Value = -0.1816402082444849
$ws.Range("F2").Value = -0.2642793407602055
$ws.Range("B3").Value = -0.4192009349369115
$ws.Range("C3").Value = -0.1408789301740885
$ws.Range("D3").Value = -0.1881079471304872
$ws.Range("E3").Value = -0.1630305981621641
$ws.Range("F3").Value = -0.1242004484768635
$ws.Range("B4").Value = -0.3135933680333336
$ws.Range("C4").Value = 0.04060011021925858
$ws.Range("D4").Value = -0.2163507239332428
$ws.Range("E4").Value = -0.1893454486309444
$ws.Range("F4").Value = -0.2261318395809019
$ws.Range("B5").Value = -0.3826290663350012
$ws.Range("C5").Value = -0.06207212619755938
$ws.Range("D5").Value = 0.08658345043665462
$ws.Range("E5").Value = 0.07752282994852688
$ws.Range("F5").Value = -0.05389568435900263
$ws.Range("B6").Value = 0.2351412773284965
$ws.Range("C6").Value = 0.3590238833875963
$ws.Range("D6").Value = 0.503602567265456
$ws.Range("E6").Value = 0.4957899563994934
$ws.Range("F6").Value = 0.3935409140535722
$ws.Range("B7").Value = 0.7417208429448655
$ws.Range("C7").Value = 0.6156611746805127
$ws.Range("D7").Value = 0.6905839940768871
$ws.Range("E7").Value = 0.7408999318973241
$ws.Range("F7").Value = 0.678643304015597
$ws.Range("B8").Value = 0.7716880594885259
$ws.Range("C8").Value = 0.7950274829098634
$ws.Range("D8").Value = 0.8945595077526429
$ws.Range("E8").Value = 0.9059220928320305
$ws.Range("F8").Value = 0.868860756151158
$ws.Range("B9").Value = 0.4069406154810908
$ws.Range("C9").Value = 0.905248353283901
$ws.Range("D9").Value = 0.9478274393755987
$ws.Range("E9").Value = 0.9478568359861684
$ws.Range("F9").Value = 0.9367287409437305
$ws.Range("B10").Value = 0.164515789745884
$ws.Range("C10").Value = 0.8940103510172752
$ws.Range("D10").Value = 0.9510893576696431
$ws.Range("E10").Value = 0.9787076460854713
$ws.Range("F10").Value = 0.9197670439975129
$ws.Range("B11").Value = -0.7347720508943657
$ws.Range("C11").Value = 0.9057016784967475
$ws.Range("D11").Value = 0.9320293786486881
$ws.Range("E11").Value = 0.9828453145286167
$ws.Range("F11").Value = 0.9385038817817638
$ws.Range("B12").Value = -0.3731218707134606
$ws.Range("C12").Value = 0.9084551194846087
$ws.Range("D12").Value = 0.9159072103908452
$ws.Range("E12").Value = 0.9802507678640193
$ws.Range("F12").Value = 0.926339891160929
$ws.Range("B13").Value = -0.9781544117523665
$ws.Range("C13").Value = 0.9013557183762829
$ws.Range("D13").Value = 0.9125994646741284
$ws.Range("E13").Value = 0.9754841628635438
$ws.Range("F13").Value = 0.9195336466068662
$ws.Range("B14").Value = -0.397251238288308
$ws.Range("C14").Value = 0.8987929643141123
$ws.Range("D14").Value = 0.9031177666895103
$ws.Range("E14").Value = 0.9716894805886346
$ws.Range("F14").Value = 0.9198067776564369
$ws.Range("B15").Value = -0.6851881314607619
$ws.Range("C15").Value = 0.8990545713409716
$ws.Range("D15").Value = 0.8909705175115198
$ws.Range("E15").Value = 0.9733799889463609
$ws.Range("F15").Value = 0.9201425940218921
$ws.Range("B16").Value = -0.6205841510403871
$ws.Range("C16").Value = 0.8994443144173788
$ws.Range("D16").Value = 0.8896356493090366
$ws.Range("E16").Value = 0.9703889827001141
$ws.Range("F16").Value = 0.9181423624008327
$ws.Range("B17").Value = -0.6328807602497957
$ws.Range("C17").Value = 0.8992193777262774
$ws.Range("D17").Value = 0.8876870587055249
$ws.Range("E17").Value = 0.9713201194571178
$ws.Range("F17").Value = 0.9177303076199722
$ws.Range("B18").Value = -0.6309414445385246
$ws.Range("C18").Value = 0.8990786983204466
$ws.Range("D18").Value = 0.8875867723379331
$ws.Range("E18").Value = 0.9711894535502141
$ws.Range("F18").Value = 0.9178531237548895
